$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Original paragraph layout (1-indexed):
#  1 Minigin things to add/do   (Title)
#  2 (empty)
#  3 BaseComponent:
#  4 (empty, bullet list numId=5)
#  5 PlayerObserver:
#  6 Fix UpdateHealthUI/GetScoreTextComponent if
#  7 (empty)
#  8 Minigin:
#  9 [list] Add fail safes to LoadScene()
# 10 [list] Improve InputManager
# 11 [list] Fix keyboard input
# 12 [list] Add prefab game objects
# 13 [list] Instead of throw std::runtime_error ... Logger class that does this
#
# Target layout:
#  7 ServiceLocator:                       (replaces empty(7) + Minigin:(8))
#  8 [list] Use template?                  (was "Add fail safes to LoadScene()")
#  9 Minigin:                              (new, plain paragraph)
# 10 [list] Add fail safes to LoadScene()  (was "Improve InputManager")
# 11 [list] Improve InputManager           (was "Fix keyboard input")
# 12 [list] Fix keyboard input             (was "Add prefab game objects")
# 13 [list] Add prefab game objects        (was "Instead of throw ...")
# 14 [list] Instead of throw ... Logger class that does this  (new list item)
# 15 (empty, new)
# 16 Implement rule of 5                   (new)
# 17 Check naming                          (new)
# 18 Cleanup Data folder                   (new)
# ------------------------------------------------------------------

# Drop the empty paragraph that used to sit before "Minigin:" and turn
# "Minigin:" into "ServiceLocator:".
$d.Paragraphs.Item(7).Range.Delete()
$d.Paragraphs.Item(7).Range.Text = "ServiceLocator:"

# The list item that used to read "Add fail safes to LoadScene()" becomes
# "Use template?" and a new plain "Minigin:" paragraph follows it.
$d.Paragraphs.Item(8).Range.Text = "Use template?"
$d.Paragraphs.Item(8).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item(9).Range.ListFormat.RemoveNumbers()
$d.Paragraphs.Item(9).Style = "Normal"
$d.Paragraphs.Item(9).Range.Text = "Minigin:"

# Shift the remaining original list items down by one slot.
$d.Paragraphs.Item(10).Range.Text = "Add fail safes to LoadScene()"
$d.Paragraphs.Item(11).Range.Text = "Improve InputManager"
$d.Paragraphs.Item(12).Range.Text = "Fix keyboard input"
$d.Paragraphs.Item(13).Range.Text = "Add prefab game objects"

# Add the new trailing list item (keeps the same list/bullet formatting).
$d.Paragraphs.Item(13).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item(14).Range.Text = "Instead of throw std::runtime_error in code, make a static Logger class that does this"

# Append the new closing block: a blank line followed by three plain
# (non-list) paragraphs.
$d.Paragraphs.Item(14).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item(15).Range.ListFormat.RemoveNumbers()
$d.Paragraphs.Item(15).Style = "Normal"

$d.Paragraphs.Item(15).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item(16).Range.ListFormat.RemoveNumbers()
$d.Paragraphs.Item(16).Style = "Normal"
$d.Paragraphs.Item(16).Range.Text = "Implement rule of 5"

$d.Paragraphs.Item(16).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item(17).Range.ListFormat.RemoveNumbers()
$d.Paragraphs.Item(17).Style = "Normal"
$d.Paragraphs.Item(17).Range.Text = "Check naming "

$d.Paragraphs.Item(17).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item(18).Range.ListFormat.RemoveNumbers()
$d.Paragraphs.Item(18).Style = "Normal"
$d.Paragraphs.Item(18).Range.Text = "Cleanup Data folder"
